$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same formatting used by the existing "CE" column to the new "CF" column
$ws.Range("CE1").Copy()
$ws.Range("CF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("CE2:CE11").Copy()
$ws.Range("CF2:CF11").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# New header / date label for the added column
$ws.Range("CF1").Value2 = "6-nov"

# New data values for the added column
$ws.Range("CF2").Value2 = 7
$ws.Range("CF3").Value2 = 10
$ws.Range("CF4").Value2 = 7
$ws.Range("CF5").Value2 = 8
$ws.Range("CF6").Value2 = 8
$ws.Range("CF7").Value2 = 10
$ws.Range("CF8").Value2 = 10
$ws.Range("CF9").Value2 = 11
$ws.Range("CF10").Value2 = 16
$ws.Range("CF11").Value2 = 0

# Update the selected cell, matching the saved view state
$ws.Range("CF12").Select()
